# Simulated Wild Card round and logged it
# Update row 3 ("R" - road game stats) on both the OFF and DEF sheets
# with the newly simulated target depth data.

$wb = $excel.ActiveWorkbook

# --- OFF sheet ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 232
$wsOff.Range("C3").Value = 160
$wsOff.Range("D3").Value = 73
$wsOff.Range("E3").Value = 31
$wsOff.Range("F3").Value = 4
$wsOff.Range("G3").Value = 5

# --- DEF sheet ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 301
$wsDef.Range("C3").Value = 235
$wsDef.Range("D3").Value = 44
$wsDef.Range("E3").Value = 21
